$wb = $excel.ActiveWorkbook

# --- Sheet1: fill rows 1..25 with repeated phone / price data ---
$ws1 = $wb.Worksheets.Item("Sheet1")
for ($r = 1; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 1).Value = "Apple iPhone 8 Plus (Gold, 64 GB)"
    $ws1.Cells.Item($r, 2).Value = "₹59,900"
}

# --- "data" sheet: add header row (Phone / Price) above the existing entry ---
$data = $wb.Worksheets.Item("data")
$data.Cells.Item(1, 1).Value = "Phone"
$data.Cells.Item(1, 2).Value = "Price"

# column widths for the "data" sheet
$data.Columns.Item(1).ColumnWidth = 44.42578125
$data.Columns.Item(2).ColumnWidth = 26.28515625

# make "data" the active / selected sheet with B1 selected
$data.Activate() | Out-Null
$data.Range("B1").Select() | Out-Null
